$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.556.65"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "1.879.70"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.026"
$ws.Range("E4").Value = "  +1.84%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.62"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.023"
$ws.Range("E6").Value = "  +1.71%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5147"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3946"
$ws.Range("E8").Value = "  +0.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08344"
$ws.Range("E9").Value = "  -0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.118"
$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.15"
$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.262"
$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.866.45"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("B14").Value = "Solana"
$ws.Range("C14").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.52"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.258"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.026"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001113"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.41"
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06778"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.72"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.024"
$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.977"
$ws.Range("E22").Value = "  -1.26%  "

$ws.Range("D23").Value = "28.603.82"
$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.268"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.078.94"
$ws.Range("E26").Value = "  -0.44%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.98"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.84"
$ws.Range("E28").Value = "  +0.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.378"
$ws.Range("E29").Value = "  -5.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.35"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1058"
$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.038"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.841"
$ws.Range("E33").Value = "  -0.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.668"
$ws.Range("E34").Value = "  +1.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02447"
$ws.Range("E35").Value = "  -0.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06527"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.161"
$ws.Range("E37").Value = "  -5.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2190"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.192"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.255"
$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6465"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.004"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.22"
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6052"
$ws.Range("E44").Value = "  -1.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.99"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  +0.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.234"
$ws.Range("E47").Value = "  -4.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.998"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.215"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.25"
$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06886"
$ws.Range("E51").Value = "  -0.34%  "
